# EPI_Poster_Deez.pptx edit script
# Applies the author's changes described in the commit:
#  - Rewrites the "METHODE" box into a long reflective paragraph and
#    resizes/un-rotates it.
#  - Moves the title box and replaces the placeholder title text.
#  - Rewrites the "EINLEITUNG" box into a new paragraph and
#    resizes/un-rotates it.
#  - Re-positions the "ERGEBNISSE" box (un-rotated).
#  - Deletes the leftover "Textfeld 16" instructions placeholder shape.
#  - Slightly re-rotates/re-positions the "Struktur Ihres Codes:" box and
#    merges its two runs into one.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "Abgerundetes Rechteck 61" (METHODE box) ---------------------
$shp = $s.Shapes.Item("Abgerundetes Rechteck 61")
$shp.Rotation = 0
$shp.Left = 692080 / 12700.0
$shp.Top = 8853443 / 12700.0
$shp.Width = 19361026 / 12700.0
$shp.Height = 4714908 / 12700.0

$tr = $shp.TextFrame.TextRange
$tr.Text = "."
$tr.Text = "Als Denkwerkzeuge wurden die Methodiken " + [char]0x201E + "On the Shoulders of Giants" + [char]0x201C + " und " + [char]0x201E + "Feedback" + [char]0x201C + " verwendet. Zur befragung haben wir Studenten und Ehemalige Studenten befragt zu diesen Themen und sie um eine Reflektion und Verbesserung ihres studiums gefragt."

# --- Shape "Abgerundetes Rechteck 1" (Title box) -------------------------
$shp = $s.Shapes.Item("Abgerundetes Rechteck 1")
$shp.Left = 334890 / 12700.0
$shp.Top = 2995527 / 12700.0

$tr = $shp.TextFrame.TextRange.Paragraphs(1, 1)
$tr.Text = "."
$tr.Text = "Stellungnahme zum Studium mit Reflektion des Gelerntem Wissen"
$tr.Font.Name = "Myriad Pro"

# --- Shape "Abgerundetes Rechteck 41" (EINLEITUNG box) -------------------
$shp = $s.Shapes.Item("Abgerundetes Rechteck 41")
$shp.Rotation = 0
$shp.Left = 620642 / 12700.0
$shp.Top = 6424551 / 12700.0
$shp.Width = 19370922 / 12700.0
$shp.Height = 2071702 / 12700.0

$tr = $shp.TextFrame.TextRange
$tr.Text = "."
$tr.Text = "In diesem Wissensgraph haben wir ein wenig Recherchiert wie der Generelle aufbau der Hochschule ist mit Wissen und Meinungen von Studenten und eine Generelle Wissensstruktur zum Thema Programmieren"

# --- Shape "Abgerundetes Rechteck 13" (ERGEBNISSE box) -------------------
$shp = $s.Shapes.Item("Abgerundetes Rechteck 13")
$shp.Rotation = 0
$shp.Left = 1049270 / 12700.0
$shp.Top = 18926201 / 12700.0

# --- Shape "Abgerundetes Rechteck 5" (Struktur Ihres Codes: box) ---------
$shp = $s.Shapes.Item("Abgerundetes Rechteck 5")
$shp.Rotation = 291993 / 60000.0
$shp.Left = 615195 / 12700.0
$shp.Top = 15537873 / 12700.0

$tr = $shp.TextFrame.TextRange
$tr.Text = "."
$tr.Text = "Struktur Ihres Codes:"

# --- Remove the leftover instructions textbox "Textfeld 16" --------------
$s.Shapes.Item("Textfeld 16").Delete()
